$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The slide's title text is already "Here is a single header" when read
# back (PowerPoint joins the runs for you), but it is still stored as nine
# separate same-formatted runs under the hood. Re-assigning the exact same
# string is treated as a no-op by the engine, so nudge the value away and
# back to force the paragraph to be rebuilt as a single run.
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "__tmp__"
$title.Text = "Here is a single header"

# Same idea for the speaker notes text box, collapsing its runs into one.
$notesRange = $s.NotesPage.Shapes.Item(2).TextFrame.TextRange
$notesRange.Text = "and here are some notes"
